$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 15) with the same shape/pattern as the
# preceding rows (row 2..14): Date in column A, numeric metrics in
# columns B..M, and the "Method" label (shared string "Noun") in column N.
$row = 15

$ws.Cells.Item($row, 1).Value = 42625.886643518519
$ws.Cells.Item($row, 2).Value = -12
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 41
$ws.Cells.Item($row, 5).Value = 33
$ws.Cells.Item($row, 6).Value = 66
$ws.Cells.Item($row, 7).Value = 13364
$ws.Cells.Item($row, 8).Value = 9708
$ws.Cells.Item($row, 9).Value = 1189
$ws.Cells.Item($row, 10).Value = 191
$ws.Cells.Item($row, 11).Value = 143
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Noun"
